# Commit: "adding rest assured dependancy"
# Adds a new "UserCredentials" worksheet with two data-provider rows,
# mirroring the layout/hyperlinks of the existing "UserDatails" sheet,
# and makes it the active tab.

$wb = $excel.ActiveWorkbook

# --- sheet1: clear the old selection/active marker (new sheet becomes active) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A1:C3").Select()
$ws1.Columns.Item(1).ColumnWidth = 20.02
$ws1.Columns.Item(2).ColumnWidth = 14.31
$ws1.Columns.Item(3).ColumnWidth = 14.74

# --- add the new sheet after the last existing sheet ---
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add($null, $last)
$ws2.Name = "UserCredentials"

# column widths for the new sheet
$ws2.Columns.Item(1).ColumnWidth = 27.31
$ws2.Columns.Item(2).ColumnWidth = 16.88
$ws2.Columns.Item(3).ColumnWidth = 16.74

# --- header row ---
$ws2.Range("A1").Value = "Test_Case_Name"
$ws2.Range("B1").Value = "UserName"
$ws2.Range("C1").Value = "Password"

# bold + fill header via a throwaway named style (keeps the style table minimal)
$headerStyle = $wb.Styles.Add("TmpHeaderStyle")
$headerStyle.Font.Bold = $true
$headerStyle.Interior.Color = 15189940
$ws2.Range("A1:C1").Style = "TmpHeaderStyle"
$headerStyle.Delete()

# --- data rows ---
$ws2.Range("A2").Value = "DataProviderWithExcel_001"
$ws2.Range("B2").Value = "testuser_1"
$ws2.Range("C2").Value = "Test@123"

$ws2.Range("A3").Value = "DataProviderWithExcel_002"
$ws2.Range("B3").Value = "testuser_1"
$ws2.Range("C3").Value = "Test@123"

# hyperlinks on the password column, matching sheet1's pattern
$ws2.Hyperlinks.Add($ws2.Range("C2"), "mailto:Test@123")
$ws2.Hyperlinks.Add($ws2.Range("C3"), "mailto:Test@123")
$ws2.Range("C2").Style = "Hyperlink"
$ws2.Range("C3").Style = "Hyperlink"

# page setup
$ws2.PageSetup.Orientation = 1

# selection + make the new sheet the active tab (matches tabSelected moving to sheet2)
$ws2.Range("C5").Select()
$ws2.Activate()
